# The post "「ベイルート・アニメイテッド」" (row 746) was removed from the sheet.
# Deleting the entire row shifts every subsequent row up by one and also
# shrinks the sheet's used range from A1:C880 to A1:C879, matching the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A746").EntireRow.Delete()
